$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Change the supplier image paths to a single shared filename
$ws.Range("D2:D5").Value = "Samsung.jpg"

# Update the selected range to reflect the new selection D2:D5
$ws.Range("D2:D5").Select()
